# FOMS order_list.xlsx - "Updated files for code logic"
# Appends a new order row (D-102 / JP / fries, fries / fries : test / false / NEW / Cash)
# to Sheet1 and refreshes the sheet's look (font, selection, page orientation) to match
# the state the workbook was saved in afterwards.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Append the new order as row 4 --------------------------------------
$newRow = $ws.Cells.Item(4, 1).EntireRow.Row

$ws.Cells.Item($newRow, 1).Value = "D-102"
$ws.Cells.Item($newRow, 2).Value = "JP"
$ws.Cells.Item($newRow, 3).Value = "fries, fries"
$ws.Cells.Item($newRow, 4).Value = "fries : test"
# Force this one to be stored as text ("false"), not as a boolean.
$ws.Cells.Item($newRow, 5).Value = "'false"
$ws.Cells.Item($newRow, 6).Value = "NEW"
$ws.Cells.Item($newRow, 7).Value = "Cash"

# --- 2. Re-apply the sheet-wide font (Calibri 11) ---------------------------
$ws.Cells.Font.Name = "Calibri"
$ws.Cells.Font.Size = 11

# --- 3. Selection state: everything selected, as left by the editing app ----
$ws.Cells.Select()

# --- 4. Page setup: explicit portrait orientation ----------------------------
$ws.PageSetup.Orientation = 1
